$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "time_taken" header in F1, matching the style of the other
# header cells (bold font, border, centered alignment - style index 1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Populate the time_taken values for each data row (no special style,
# same as the other plain data cells)
$ws.Range("F2").Value = "2021-10-05 10:51:56.384158"
$ws.Range("F3").Value = "2021-10-05 10:51:56.384170"
$ws.Range("F4").Value = "2021-10-05 10:51:56.384174"
$ws.Range("F5").Value = "2021-10-05 10:51:56.384177"
$ws.Range("F6").Value = "2021-10-05 10:51:56.384180"
$ws.Range("F7").Value = "2021-10-05 10:51:56.384183"
$ws.Range("F8").Value = "2021-10-05 10:51:56.384186"
$ws.Range("F9").Value = "2021-10-05 10:51:56.384189"
$ws.Range("F10").Value = "2021-10-05 10:51:56.384192"
$ws.Range("F11").Value = "2021-10-05 10:51:56.384196"
$ws.Range("F12").Value = "2021-10-05 10:51:56.384199"
$ws.Range("F13").Value = "2021-10-05 10:51:56.384201"
$ws.Range("F14").Value = "2021-10-05 10:51:56.384204"
$ws.Range("F15").Value = "2021-10-05 10:51:56.384207"
$ws.Range("F16").Value = "2021-10-05 10:51:56.384210"
$ws.Range("F17").Value = "2021-10-05 10:51:56.384213"
